# Leave Card update: add 2024 accrual entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item(1)

# --- Fill in the EARNED values that were posted for the first few
#     2023-2024 VL(1-0-0) periods (rows 22-25) ---
$ws.Range("C22").Value = 1.25
$ws.Range("C23").Value = 1.25
$ws.Range("C24").Value = 1.25
$ws.Range("C25").Value = 1.25

# --- Row 26 (12/21/2023 period): mark a 3-day VL ---
$ws.Range("B26").Value = "VL(3-0-0)"
$ws.Range("D26").Value = 3

# --- Insert a new "2024" year-divider row above the old row 27,
#     shifting every following row down by one (table grows by a row) ---
$ws.Rows("27:27").Insert()

# Copy the formatting from the existing "2023" divider row (row 14) so the
# new divider row matches the established look.
$ws.Range("A14:K14").Copy()
$ws.Range("A27:K27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A27").Value = "2024"
$ws.Range("G27").Formula = "=IF(ISBLANK([@EARNED]),"",[@EARNED])"

# Resize the table to include the newly inserted row.
$newTableRange = $ws.Range($lo.Range.Cells.Item(1, 1), $ws.Range("K133"))
$lo.Resize($newTableRange)

$wb.Save()
